# Commit: "Actividad practica laravel 22/11/2022"
#
# The paragraph "Por último solo nos queda indicar una redirección..."
# had two stray grammar-check markers (<w:proofErr w:type="gramStart"/> and
# <w:proofErr w:type="gramEnd"/>) bracketing the word "último". The edit
# drops those proofing markers and adds a comma right after "último"
# ("Por último," instead of "Por último").

$d = $word.ActiveDocument

# This exact sentence is unique in the document, so Find will land on the
# right paragraph (there is another, unrelated "último" later in the doc).
$target = "Por último solo nos queda indicar una redirección la cual, después de añadir a la base de datos se ejecuta de la siguiente forma "

$rng = $d.Content
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Output "WARNING: target sentence not found; document left unchanged"
}
else {
    # Remove the old text (this also drops the gramStart/gramEnd proofErr
    # markers that sat inside the deleted range) ...
    $rng.Delete()

    # ... and retype it as three plain runs ("Por ", "último,", " solo nos
    # queda ... forma ") with no proofing markers, matching the target markup.
    $insertionPoint = $d.Range($rng.Start, $rng.Start)
    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Por </w:t></w:r>
            <w:r><w:t>último,</w:t></w:r>
            <w:r><w:t xml:space="preserve"> solo nos queda indicar una redirección la cual, después de añadir a la base de datos se ejecuta de la siguiente forma </w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
    $insertionPoint.InsertXML($xml)

    Write-Output "Replaced 'Por [proofErr]ultimo[proofErr]' with 'Por ultimo,' (proofErr markers removed)"
}
